$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "58.995.38"
Set-TextValue "E2" "  +0.91%  "

# Row 3
Set-TextValue "D3" "2.640.07"
Set-TextValue "E3" "  +4.36%  "

# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.12%  "

# Row 5
Set-TextValue "D5" "517.28"
Set-TextValue "E5" "  +2.14%  "

# Row 6
Set-TextValue "D6" "144.31"
Set-TextValue "E6" "  +0.58%  "

# Row 7
Set-TextValue "D7" "0.997"
Set-TextValue "E7" "  -0.27%  "

# Row 8
Set-TextValue "D8" "0.569"
Set-TextValue "E8" "  +1.56%  "

# Row 9
Set-TextValue "D9" "2.662.56"
Set-TextValue "E9" "  +5.07%  "

# Row 10
Set-TextValue "D10" "6.25"
Set-TextValue "E10" "  +1.40%  "

# Row 11
Set-TextValue "D11" "0.105"
Set-TextValue "E11" "  +3.75%  "

# Row 12
Set-TextValue "D12" "0.338"
Set-TextValue "E12" "  +2.23%  "

# Row 13
Set-TextValue "E13" "  -1.67%  "

# Row 14
Set-TextValue "D14" "3.110.21"
Set-TextValue "E14" "  +4.51%  "

# Row 15
Set-TextValue "D15" "58.941.51"
Set-TextValue "E15" "  +0.77%  "

# Row 16
Set-TextValue "D16" "21.03"
Set-TextValue "E16" "  +2.05%  "

# Row 17
Set-TextValue "D17" "0.0000137"
Set-TextValue "E17" "  +1.85%  "

# Row 18
Set-TextValue "D18" "2.666.13"
Set-TextValue "E18" "  +5.16%  "

# Row 19
Set-TextValue "D19" "351.96"
Set-TextValue "E19" "  +5.25%  "

# Row 20
Set-TextValue "D20" "4.54"
Set-TextValue "E20" "  +0.50%  "

# Row 21
Set-TextValue "D21" "10.39"
Set-TextValue "E21" "  +3.49%  "

# Row 22
Set-TextValue "D22" "6.20"
Set-TextValue "E22" "  +4.39%  "

# Row 23
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.24%  "

# Row 24
Set-TextValue "D24" "61.88"
Set-TextValue "E24" "  +3.03%  "

# Row 25
Set-TextValue "D25" "0.421"
Set-TextValue "E25" "  +3.26%  "

# Row 26
Set-TextValue "D26" "0.997"
Set-TextValue "E26" "  -0.20%  "

# Row 27
Set-TextValue "E27" "  +0.99%  "

# Row 28
Set-TextValue "D28" "0.0₃0808"
Set-TextValue "E28" "  +3.13%  "

# Row 29
Set-TextValue "D29" "7.19"
Set-TextValue "E29" "  +3.87%  "

# Row 30
Set-TextValue "D30" "0.998"
Set-TextValue "E30" "  -0.16%  "

# Row 31
Set-TextValue "D31" "6.33"
Set-TextValue "E31" "  +8.82%  "

# Row 32
Set-TextValue "D32" "19.05"
Set-TextValue "E32" "  +3.11%  "

# Row 33
Set-TextValue "D33" "1.58"
Set-TextValue "E33" "  +3.11%  "

# Row 34
Set-TextValue "D34" "150.36"
Set-TextValue "E34" "  +0.43%  "

# Row 35
Set-TextValue "D35" "0.970"
Set-TextValue "E35" "  +3.82%  "

# Row 36
Set-TextValue "D36" "4.02"

# Row 37
Set-TextValue "D37" "1.14"
Set-TextValue "E37" "  +2.84%  "

# Row 38
Set-TextValue "D38" "36.60"
Set-TextValue "E38" "  +1.57%  "

# Row 39
Set-TextValue "D39" "0.844"
Set-TextValue "E39" "  +2.50%  "

# Row 40
Set-TextValue "D40" "3.74"
Set-TextValue "E40" "  +6.33%  "

# Row 41
Set-TextValue "D41" "1.41"
Set-TextValue "E41" "  +1.12%  "

# Row 42
Set-TextValue "D42" "281.40"
Set-TextValue "E42" "  -0.27%  "

# Row 43
Set-TextValue "D43" "0.613"
Set-TextValue "E43" "  +1.98%  "

# Row 44
Set-TextValue "D44" "0.0989"
Set-TextValue "E44" "  -0.45%  "

# Row 45
Set-TextValue "B45" "FirstDigitalUSD"
Set-TextValue "C45" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D45" "0.994"
Set-TextValue "E45" "  -0.43%  "

# Row 46
Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "19.63"
Set-TextValue "E46" "  +5.49%  "

# Row 47
Set-TextValue "D47" "0.0532"
Set-TextValue "E47" "  -0.11%  "

# Row 48
Set-TextValue "D48" "0.0231"
Set-TextValue "E48" "  +2.24%  "

# Row 49
Set-TextValue "B49" "WhiteBITCoin"
Set-TextValue "C49" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D49" "10.28"
Set-TextValue "E49" "  -0.20%  "

# Row 50
Set-TextValue "B50" "Maker"
Set-TextValue "C50" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "1.990.75"
Set-TextValue "E50" "  +5.32%  "

# Row 51
Set-TextValue "B51" "RenderToken"
Set-TextValue "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "4.68"
Set-TextValue "E51" "  +3.60%  "
